$d = $word.ActiveDocument

# The three paragraphs:
#   P2: "The Department of NOVA is responsible for the planning and overseeing of the annual NOVA"
#   P3: "banquet meant to celebrate the students and staff of Florida Polytechnic University and their"
#   P4: "accomplishments."
# get merged into a single paragraph, with the two paragraph breaks between
# them replaced by single space characters, each kept as its own run.
#
# Work from the end of the document towards the start so earlier character
# offsets stay valid while later ones are being edited.

$p2 = $d.Paragraphs(2)
$p3 = $d.Paragraphs(3)

$junction2 = $p3.Range.End - 1   # position of the paragraph mark ending P3
$junction1 = $p2.Range.End - 1   # position of the paragraph mark ending P2

# --- second junction: "...and their" | "accomplishments." ---
$mark2 = $d.Range($junction2, $junction2 + 1)
$mark2.Delete()

$spacePt2 = $d.Range($junction2, $junction2)
$spacePt2.InsertBefore(" ")
$spaceRun2 = $d.Range($junction2, $junction2 + 1)
$spaceRun2.Font.Bold = $true
$spaceRun2.Font.Bold = $false

# --- first junction: "...annual NOVA" | "banquet..." ---
$mark1 = $d.Range($junction1, $junction1 + 1)
$mark1.Delete()

$spacePt1 = $d.Range($junction1, $junction1)
$spacePt1.InsertBefore(" ")
$spaceRun1 = $d.Range($junction1, $junction1 + 1)
$spaceRun1.Font.Bold = $true
$spaceRun1.Font.Bold = $false

$d.Save()
